$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 3042
$ws.Range("I34").Value = 3042
$ws.Range("K34").Value = 3042
$ws.Range("M34").Value = -2839
$ws.Range("H36").Value = 3042
$ws.Range("I36").Value = 3042
$ws.Range("K36").Value = 3042
$ws.Range("M36").Value = -2327
$ws.Range("H64").Value = 7500
$ws.Range("J64").Value = 12000
$ws.Range("L64").Value = 12000
$ws.Range("N64").Value = -12496
$ws.Range("H67").Value = 7500
$ws.Range("J67").Value = 12000
$ws.Range("L67").Value = 12000
$ws.Range("N67").Value = -13716
$ws.Range("H98").Value = 1544.6538
$ws.Range("I98").Value = 1329.7
$ws.Range("J98").Value = 2261.1667
$ws.Range("K98").Value = 1329.7
$ws.Range("L98").Value = 2261.1667
$ws.Range("M98").Value = 168.3
$ws.Range("N98").Value = -5257.1667
$ws.Range("H122").Value = 1544.6538
$ws.Range("I122").Value = 1329.7
$ws.Range("J122").Value = 2261.1667
$ws.Range("K122").Value = 3989.1
$ws.Range("L122").Value = 6783.500100000001
$ws.Range("M122").Value = -1539.1
$ws.Range("N122").Value = -11683.5001
$ws.Range("H137").Value = 4333170.5
$ws.Range("I137").Value = 5000
$ws.Range("K137").Value = 15000
$ws.Range("M137").Value = -12450

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2549.7693
$ws.Range("I45").Value = 1340.1
$ws.Range("K45").Value = 1340.1
$ws.Range("M45").Value = -963.0999999999999
$ws.Range("H46").Value = 18150.637
$ws.Range("J46").Value = 18975.9
$ws.Range("L46").Value = 18975.9
$ws.Range("N46").Value = -19613.9
$ws.Range("H101").Value = 122966
$ws.Range("J101").Value = 122966
$ws.Range("L101").Value = 122966
$ws.Range("N101").Value = -129456
$ws.Range("H124").Value = 34613.5
$ws.Range("J124").Value = 34613.5
$ws.Range("L124").Value = 34613.5
$ws.Range("N124").Value = -44433.5
$ws.Range("H125").Value = 63611
$ws.Range("J125").Value = 63611
$ws.Range("L125").Value = 63611
$ws.Range("N125").Value = -73451

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 175000
$ws.Range("J132").Value = 175000
$ws.Range("L132").Value = 175000
$ws.Range("N132").Value = -185120
$ws.Range("H134").Value = 12502864
$ws.Range("I134").Value = 2982.4
$ws.Range("K134").Value = 8947.200000000001
$ws.Range("M134").Value = -6412.200000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2643.1353
$ws.Range("I107").Value = 2456.6785
$ws.Range("J107").Value = 3223.2222
$ws.Range("K107").Value = 2456.6785
$ws.Range("L107").Value = 3223.2222
$ws.Range("M107").Value = -536.6785
$ws.Range("N107").Value = -7063.2222
$ws.Range("H135").Value = 129999.25
$ws.Range("J135").Value = 129999.25
$ws.Range("L135").Value = 129999.25
$ws.Range("N135").Value = -140139.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 19236618
$ws.Range("I129").Value = 25003072
$ws.Range("J129").Value = 15107.333
$ws.Range("K129").Value = 75009216
$ws.Range("L129").Value = 45321.999
$ws.Range("M129").Value = -75004216
$ws.Range("N129").Value = -55321.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 2428.1428
$ws.Range("I107").Value = 2400
$ws.Range("J107").Value = 2465.6667
$ws.Range("K107").Value = 2400
$ws.Range("L107").Value = 2465.6667
$ws.Range("M107").Value = -480
$ws.Range("N107").Value = -6305.6667
$ws.Range("H133").Value = 125779
$ws.Range("J133").Value = 125779
$ws.Range("L133").Value = 125779
$ws.Range("N133").Value = -135899
$ws.Range("H136").Value = 78438
$ws.Range("J136").Value = 78438
$ws.Range("L136").Value = 235314
$ws.Range("N136").Value = -240414

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6516.8335
$ws.Range("I40").Value = 6520.2
$ws.Range("K40").Value = 6520.2
$ws.Range("M40").Value = -6384.2
$ws.Range("H61").Value = 9064.429
$ws.Range("I61").Value = 2221.5
$ws.Range("J61").Value = 11801.6
$ws.Range("K61").Value = 2221.5
$ws.Range("L61").Value = 11801.6
$ws.Range("M61").Value = -2019.5
$ws.Range("N61").Value = -12205.6
$ws.Range("H68").Value = 4631970.5
$ws.Range("I68").Value = 6946272
$ws.Range("J68").Value = 3367.6667
$ws.Range("K68").Value = 6946272
$ws.Range("L68").Value = 3367.6667
$ws.Range("M68").Value = -6945523
$ws.Range("N68").Value = -4865.6667
$ws.Range("H71").Value = 4631970.5
$ws.Range("I71").Value = 6946272
$ws.Range("J71").Value = 3367.6667
$ws.Range("K71").Value = 34731360
$ws.Range("L71").Value = 16838.3335
$ws.Range("M71").Value = -34727616
$ws.Range("N71").Value = -24326.3335
$ws.Range("H82").Value = 6328
$ws.Range("I82").Value = 4380
$ws.Range("K82").Value = 4380
$ws.Range("M82").Value = -4019
$ws.Range("H85").Value = 6328
$ws.Range("I85").Value = 4380
$ws.Range("K85").Value = 4380
$ws.Range("M85").Value = -3132
$ws.Range("H113").Value = 9064.429
$ws.Range("I113").Value = 2221.5
$ws.Range("J113").Value = 11801.6
$ws.Range("K113").Value = 2221.5
$ws.Range("L113").Value = 11801.6
$ws.Range("M113").Value = -51.5
$ws.Range("N113").Value = -16141.6
$ws.Range("H136").Value = 3304
$ws.Range("I136").Value = 2139.4666
$ws.Range("J136").Value = 5487.5
$ws.Range("K136").Value = 6418.399800000001
$ws.Range("L136").Value = 16462.5
$ws.Range("M136").Value = -3868.399800000001
$ws.Range("N136").Value = -21562.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 29992.666
$ws.Range("J41").Value = 29992.5
$ws.Range("L41").Value = 29992.5
$ws.Range("N41").Value = -30772.5
$ws.Range("H62").Value = 13745
$ws.Range("I62").Value = 4990
$ws.Range("J62").Value = 16663.334
$ws.Range("K62").Value = 4990
$ws.Range("L62").Value = 16663.334
$ws.Range("M62").Value = -4366
$ws.Range("N62").Value = -17911.334
$ws.Range("H65").Value = 13745
$ws.Range("I65").Value = 4990
$ws.Range("J65").Value = 16663.334
$ws.Range("K65").Value = 24950
$ws.Range("L65").Value = 83316.67
$ws.Range("M65").Value = -21830
$ws.Range("N65").Value = -89556.67
$ws.Range("H74").Value = 15071.714
$ws.Range("J74").Value = 14008.5
$ws.Range("L74").Value = 14008.5
$ws.Range("N74").Value = -15880.5
$ws.Range("H77").Value = 15071.714
$ws.Range("J77").Value = 14008.5
$ws.Range("L77").Value = 42025.5
$ws.Range("N77").Value = -51385.5
$ws.Range("H81").Value = 1668.3125
$ws.Range("I81").Value = 1512.8667
$ws.Range("K81").Value = 3025.7334
$ws.Range("M81").Value = -1964.7334
$ws.Range("H84").Value = 1668.3125
$ws.Range("I84").Value = 1512.8667
$ws.Range("K84").Value = 15128.667
$ws.Range("M84").Value = -9824.667000000001
$ws.Range("H100").Value = 1385.5385
$ws.Range("I100").Value = 1359.3334
$ws.Range("K100").Value = 2718.6668
$ws.Range("M100").Value = -2177.6668
$ws.Range("H122").Value = 2305.423
$ws.Range("I122").Value = 1853.7646
$ws.Range("J122").Value = 3158.5557
$ws.Range("K122").Value = 5561.293799999999
$ws.Range("L122").Value = 9475.667099999999
$ws.Range("M122").Value = -3111.293799999999
$ws.Range("N122").Value = -14375.6671
$ws.Range("H136").Value = 287039.5
$ws.Range("I136").Value = 1440.0968
$ws.Range("K136").Value = 4320.2904
$ws.Range("M136").Value = -1770.2904
